# Update contestant record day groups: append 4 new Group rows
# (GRP002..GRP005) to the "Groups" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Groups")

$newGroups = @(
    @("8a692ea1-5300-494c-8074-19f479c5284d", "GRP002"),
    @("6f3cbd80-aae5-463e-836b-f1ecf0c504a2", "GRP003"),
    @("3e3161fd-95f1-4ce9-9810-0ab1d377e2b2", "GRP004"),
    @("5b77330c-45dc-4177-b704-50e8d2e75d75", "GRP005")
)

$startRow = 3
for ($i = 0; $i -lt $newGroups.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $newGroups[$i][0]
    $ws.Cells.Item($row, 2).Value = $newGroups[$i][1]
}
